# Auto-generated edit script
# Applies the numeric value updates described by the commit diff
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR sheets.

$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 70
$ws.Range("H70").Value = 568799.3
$ws.Range("I70").Value = 3402269.2
$ws.Range("J70").Value = 2105.3333
$ws.Range("K70").Value = 10206807.6
$ws.Range("L70").Value = 6315.999899999999
$ws.Range("M70").Value = -10206537.6
$ws.Range("N70").Value = -6855.999899999999
# row 73
$ws.Range("H73").Value = 568799.3
$ws.Range("I73").Value = 3402269.2
$ws.Range("J73").Value = 2105.3333
$ws.Range("K73").Value = 10206807.6
$ws.Range("L73").Value = 6315.999899999999
$ws.Range("M73").Value = -10205871.6
$ws.Range("N73").Value = -8187.999899999999
# row 111
$ws.Range("H111").Value = 56740.5
$ws.Range("I111").Value = 1899
$ws.Range("K111").Value = 5697
$ws.Range("M111").Value = -2630
# row 116
$ws.Range("H116").Value = 58446500
$ws.Range("I116").Value = 50204000
$ws.Range("K116").Value = 50204000
$ws.Range("M116").Value = -50200558
# row 131
$ws.Range("H131").Value = 7666.75
$ws.Range("I131").Value = 699.9091
$ws.Range("K131").Value = 2099.7273
$ws.Range("M131").Value = 2940.2727
# row 132
$ws.Range("H132").Value = 2868.88
$ws.Range("I132").Value = 2525.068
$ws.Range("K132").Value = 7575.204000000001
$ws.Range("M132").Value = -5045.204000000001

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 45
$ws.Range("H45").Value = 2539.6
$ws.Range("I45").Value = 2424.5
$ws.Range("K45").Value = 2424.5
$ws.Range("M45").Value = -2047.5
# row 56
$ws.Range("H56").Value = 0
$ws.Range("I56").Value = 0
$ws.Range("K56").Value = 0
$ws.Range("M56").ClearContents()
# row 61
$ws.Range("H61").Value = 17545996
$ws.Range("I61").Value = 22223728
$ws.Range("K61").Value = 22223728
$ws.Range("M61").Value = -22223516
# row 63
$ws.Range("H63").Value = 150016640
$ws.Range("J63").Value = 66697900
$ws.Range("L63").Value = 66697900
$ws.Range("N63").Value = -66699272
# row 66
$ws.Range("H66").Value = 150016640
$ws.Range("J66").Value = 66697900
$ws.Range("L66").Value = 333489500
$ws.Range("N66").Value = -333496364
# row 74
$ws.Range("H74").Value = 1990.76
$ws.Range("I74").Value = 1687.2632
$ws.Range("K74").Value = 1687.2632
$ws.Range("M74").Value = -813.2632000000001
# row 77
$ws.Range("H77").Value = 1990.76
$ws.Range("I77").Value = 1687.2632
$ws.Range("K77").Value = 8436.316000000001
$ws.Range("M77").Value = -4068.316000000001
# row 132
$ws.Range("H132").Value = 17242710
$ws.Range("I132").Value = 18869128
$ws.Range("K132").Value = 56607384
$ws.Range("M132").Value = -56604854
# row 136
$ws.Range("H136").Value = 17545996
$ws.Range("I136").Value = 22223728
$ws.Range("K136").Value = 66671184
$ws.Range("M136").Value = -66668634
# row 139
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 94
$ws.Range("H94").Value = 2548.9167
$ws.Range("J94").Value = 4997.5
$ws.Range("L94").Value = 4997.5
$ws.Range("N94").Value = -5899.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 105
$ws.Range("H105").Value = 1683.4286
$ws.Range("I105").Value = 1630.6666
$ws.Range("K105").Value = 1630.6666
$ws.Range("M105").Value = 116.3334
# row 107
$ws.Range("H107").Value = 2170.76
$ws.Range("I107").Value = 1769.9048
$ws.Range("J107").Value = 4275.25
$ws.Range("K107").Value = 1769.9048
$ws.Range("L107").Value = 4275.25
$ws.Range("M107").Value = 150.0952
$ws.Range("N107").Value = -8115.25
# row 129
$ws.Range("H129").Value = 58890
$ws.Range("J129").Value = 58890
$ws.Range("L129").Value = 58890
$ws.Range("N129").Value = -68890

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 122
$ws.Range("H122").Value = 821.5714
$ws.Range("I122").Value = 714.25
$ws.Range("J122").Value = 964.6667
$ws.Range("K122").Value = 6428.25
$ws.Range("L122").Value = 8682.0003
$ws.Range("M122").Value = -3978.25
$ws.Range("N122").Value = -13582.0003
# row 131
$ws.Range("H131").Value = 8510.272000000001
$ws.Range("I131").Value = 1199.8
$ws.Range("J131").Value = 9815.714
$ws.Range("K131").Value = 3599.4
$ws.Range("L131").Value = 29447.142
$ws.Range("M131").Value = 1440.6
$ws.Range("N131").Value = -39527.142
# row 132
$ws.Range("H132").Value = 1809.091
$ws.Range("I132").Value = 1500
$ws.Range("J132").Value = 3200
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 28800
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -33860

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 80
$ws.Range("H80").Value = 3235.75
$ws.Range("I80").Value = 3949.5
$ws.Range("J80").Value = 2997.8333
$ws.Range("K80").Value = 3949.5
$ws.Range("L80").Value = 2997.8333
$ws.Range("M80").Value = -2951.5
$ws.Range("N80").Value = -4993.8333
# row 83
$ws.Range("H83").Value = 3235.75
$ws.Range("I83").Value = 3949.5
$ws.Range("J83").Value = 2997.8333
$ws.Range("K83").Value = 19747.5
$ws.Range("L83").Value = 14989.1665
$ws.Range("M83").Value = -14755.5
$ws.Range("N83").Value = -24973.1665
# row 113
$ws.Range("H113").Value = 2725
$ws.Range("J113").Value = 2971.4285
$ws.Range("L113").Value = 2971.4285
$ws.Range("N113").Value = -7311.4285
# row 132
$ws.Range("H132").Value = 6623.9287
$ws.Range("I132").Value = 8600
$ws.Range("J132").Value = 5526.1113
$ws.Range("K132").Value = 25800
$ws.Range("L132").Value = 16578.3339
$ws.Range("M132").Value = -23270
$ws.Range("N132").Value = -21638.3339
# row 134
$ws.Range("H134").Value = 109831.5
$ws.Range("J134").Value = 109831.5
$ws.Range("L134").Value = 329494.5
$ws.Range("N134").Value = -334564.5

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 40
$ws.Range("H40").Value = 3222
$ws.Range("I40").Value = 2796.4167
$ws.Range("J40").Value = 4498.75
$ws.Range("K40").Value = 2796.4167
$ws.Range("L40").Value = 4498.75
$ws.Range("M40").Value = -2660.4167
$ws.Range("N40").Value = -4770.75
# row 46
$ws.Range("H46").Value = 3125.1333
$ws.Range("I46").Value = 2232.8333
$ws.Range("J46").Value = 3720
$ws.Range("K46").Value = 2232.8333
$ws.Range("L46").Value = 3720
$ws.Range("M46").Value = -2044.8333
$ws.Range("N46").Value = -4096
# row 50
$ws.Range("H50").Value = 25000
$ws.Range("I50").Value = 25000
$ws.Range("K50").Value = 25000
$ws.Range("M50").Value = -24363
# row 132
$ws.Range("H132").Value = 4363.5386
$ws.Range("I132").Value = 3481.8948
$ws.Range("K132").Value = 10445.6844
$ws.Range("M132").Value = -7915.6844

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 94
$ws.Range("H94").Value = 25999.334
$ws.Range("J94").Value = 25999.334
$ws.Range("L94").Value = 25999.334
$ws.Range("N94").Value = -27801.334
# row 96
$ws.Range("H96").Value = 2377
$ws.Range("J96").Value = 2377
$ws.Range("L96").Value = 2377
$ws.Range("N96").Value = -5123
# row 106
$ws.Range("H106").Value = 21998.5
$ws.Range("I106").Value = 24999
$ws.Range("J106").Value = 20998.334
$ws.Range("K106").Value = 24999
$ws.Range("L106").Value = 20998.334
$ws.Range("M106").Value = -23737
$ws.Range("N106").Value = -23522.334
# row 122
$ws.Range("H122").Value = 2364.6843
$ws.Range("I122").Value = 2266.3076
$ws.Range("J122").Value = 2577.8333
$ws.Range("K122").Value = 6798.9228
$ws.Range("L122").Value = 7733.499899999999
$ws.Range("M122").Value = -4348.9228
$ws.Range("N122").Value = -12633.4999
# row 132
$ws.Range("H132").Value = 4902.55
$ws.Range("J132").Value = 3499.5715
$ws.Range("L132").Value = 10498.7145
$ws.Range("N132").Value = -15558.7145

